$wb = $excel.ActiveWorkbook

# --- CollectionEvents sheet: merge ageMin/ageMax columns into a single "ageGroups" column ---
$ws = $wb.Worksheets.Item("CollectionEvents")

# Remove the "ageMax" column (E) - "ageMin" (D) becomes "ageGroups", "subcohorts" shifts left to E
$ws.Columns.Item(5).Delete()
$ws.Cells.Item(1, 4).Value = "ageGroups"

# Make CollectionEvents the active sheet with the given selection
$ws.Activate()
$ws.Range("C4").Select()
